$d = $word.ActiveDocument

function Get-ParagraphByExactText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq $text) {
            return $p
        }
    }
    return $null
}

function Get-ParagraphByPrefix($doc, [string]$prefix) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Creer un dossier ou le serveur se trouvera (mkdir Copernicus)"
#    -> "Creer l'environnement de travail : mkdir -p Copernicus/serveur"
# ---------------------------------------------------------------------------
$p1 = Get-ParagraphByExactText $d "Créer un dossier où le serveur se trouvera (mkdir Copernicus)`r"
if ($p1 -ne $null) {
    $p1.Range.Text = "Créer l’environnement de travail : mkdir -p Copernicus/serveur"
}

# ---------------------------------------------------------------------------
# 2) "Se placer dans ce dossier" -> "Se placer dans le dossier du server : cd Copernicus/serveur"
#    (careful: another paragraph starts with the same words later on)
# ---------------------------------------------------------------------------
$p2 = Get-ParagraphByExactText $d "Se placer dans ce dossier`r"
if ($p2 -ne $null) {
    $p2.Range.Text = "Se placer dans le dossier du server : cd Copernicus/serveur"
}

# ---------------------------------------------------------------------------
# 3) "Placer dans le dossier le fichier app.py fourni (" ->
#    "Placer dans le dossier "Serveur" le fichier app.py fourni ("
#    This paragraph also contains a hyperlink run after this text, so we only
#    rewrite the leading run's text via a sub-range, leaving the hyperlink and
#    closing parenthesis run untouched.
# ---------------------------------------------------------------------------
$oldStart3 = "Placer dans le dossier le fichier app.py fourni ("
$p3 = Get-ParagraphByPrefix $d $oldStart3
if ($p3 -ne $null) {
    $start3 = $p3.Range.Start
    $sub3 = $d.Range($start3, $start3 + $oldStart3.Length)
    $sub3.Text = "Placer dans le dossier “Serveur” le fichier app.py fourni ("
}

# ---------------------------------------------------------------------------
# 4) Insert four new paragraphs after the FIRST "pip install cdsapi [Version
#    testée 0.5.1]" paragraph, before the "Vous pouvez maintenant lancer le
#    serveur" (numId=6) paragraph.
# ---------------------------------------------------------------------------
$cds1 = Get-ParagraphByExactText $d "pip install cdsapi [Version testée 0.5.1]`r"
if ($cds1 -ne $null) {
    $anchor = $cds1
    $texts4 = @(
        "pip install xarray [Version testée 2022.3.0]",
        "pip install rasterio==1.3a3 [Version forcée 1.3a3]",
        "pip install rioaxarray [Version testée 0.11.1]",
        "pip install netcdf4 [Version testée 1.5.8]"
    )
    foreach ($t in $texts4) {
        $anchor.Range.InsertParagraphAfter()
        $newPara = $anchor.Next()
        $newPara.Range.Text = $t
        $anchor = $newPara
    }
}

# ---------------------------------------------------------------------------
# 5) Insert a new paragraph after "Créer un dossier dans "Documents" nommé
#    Copernicus via l'explorateur Windows", before "Placer dans ce dossier le
#    fichier app.py fourni (...)".
# ---------------------------------------------------------------------------
$docFolder = Get-ParagraphByExactText $d "Créer un dossier dans “Documents” nommé Copernicus via l’explorateur Windows`r"
if ($docFolder -ne $null) {
    $docFolder.Range.InsertParagraphAfter()
    $newSub = $docFolder.Next()
    $newSub.Range.Text = "Dans ce même dossier, créer un sous dossier nommé “serveur”"
    $newSub.Range.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 6) "cd Documents\Copernicus_Flask" -> "cd Documents\Copernicus\serveur"
# ---------------------------------------------------------------------------
$p6 = Get-ParagraphByExactText $d "cd Documents\Copernicus_Flask`r"
if ($p6 -ne $null) {
    $p6.Range.Text = "cd Documents\Copernicus\serveur"
}

# ---------------------------------------------------------------------------
# 7) Insert four new paragraphs after the SECOND "pip install cdsapi [Version
#    testée 0.5.1]" paragraph, before "Vous pouvez maintenant lancer le
#    serveur : python app.py" (numId=5) paragraph.
# ---------------------------------------------------------------------------
$cds2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "pip install cdsapi [Version testée 0.5.1]`r") {
        $cds2 = $p
    }
}
if ($cds2 -ne $null) {
    $anchor2 = $cds2
    $texts7 = @(
        "pip install xarray [Version testée 2022.3.0]",
        "pip install rasterio==1.3a3 [Version forcée 1.3a3]",
        "pip install rioaxarray [Version testée 0.11.1]",
        "pip install netcdf4 [Version testée 1.5.8]"
    )
    foreach ($t in $texts7) {
        $anchor2.Range.InsertParagraphAfter()
        $newPara2 = $anchor2.Next()
        $newPara2.Range.Text = $t
        $anchor2 = $newPara2
    }
}

Write-Output "done"
